$d = $word.ActiveDocument

# --- Change 1: paragraph 3 ("1. Worldwide, Theater...") ---
# "with " moves from before "Kickstarter campaign category" to after it.
$rng = $d.Paragraphs(3).Range
$null = $rng.Find.Execute(
    "is the most popular with Kickstarter campaign category the most total",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "is the most popular Kickstarter campaign category with the most total",
    2)

# --- Change 2: paragraph 5 ("3. Successful campaigns trend downward...") ---
# Reorders / rewrites the closing sentence about December.
$rng = $d.Paragraphs(5).Range
$null = $rng.Find.Execute(
    "December is the only month when of the campaigns started, the Failed campaigns outnumbered Successful campaigns. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Of the campaigns started in a given month, December is the only month when, the Failed campaigns outnumbered Successful campaigns. ",
    2)

# --- Change 3: paragraph 7, part A ("The first limitation...") ---
# Inserts a new sentence about the 4,000-project sample.
$rng = $d.Paragraphs(7).Range
$null = $rng.Find.Execute(
    "is not all 300,000 projects are being analyzed. Another limitation",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "is not all 300,000 projects are being analyzed. Given a different set of 4,000 projects could yield different results and conclusions. Another limitation",
    2)

# --- Change 4: paragraph 7, part B ("...Not all countries and regions...") ---
# The old text was split "N" | "ot all countries...campaigns.  " across a
# run boundary (with the old "_GoBack" bookmark sitting between them). The
# find string spans both runs so they collapse into one "Not all..." run;
# removes the old sentence-ending period and appends a new clause + trailing
# spaces. (The "_GoBack" bookmark that used to live here gets re-added
# elsewhere, below.)
$rng = $d.Paragraphs(7).Range
$null = $rng.Find.Execute(
    "Not all countries and regions are represented. Given Kickstarter is a crowdfunding platform, disposable income would be important to participate in the campaigns.  ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Not all countries and regions are represented. Given Kickstarter is a crowdfunding platform, disposable income would be important to participate in the campaigns so not all regions would have the same opportunity to fund successful campaigns.    ",
    2)

# --- Change 5: paragraph 10 ("There are a multitude of other tables/graphs...") ---
# Rewrites the "trending up or down..." sentence and fixes "and table" -> "a table".
$rng = $d.Paragraphs(10).Range
$null = $rng.Find.Execute(
    "in a given campaign status. You could create and table",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "over the years and to see which campaigns are most popular now. You could create a table",
    2)

# Re-seat the "_GoBack" bookmark: it now belongs right after
# "...most popular now. You could create a" (and before " table and graph...").
$findRng = $d.Paragraphs(10).Range.Duplicate
$null = $findRng.Find.Execute(
    "most popular now. You could create a",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmPos = $findRng.End
$bmRng = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRng)
